$wb = $excel.ActiveWorkbook

# Mapping of sheet name -> cell -> new value for the "想去人数" (want-to-go count) column
$updates = @{
    "展览"   = @{ "F2" = 371; "F3" = 797; "F5" = 891; "F6" = 2174; "F7" = 192 }
    "全部类型" = @{ "F2" = 371; "F3" = 797; "F7" = 891; "F8" = 2174; "F10" = 192 }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cellMap = $updates[$sheetName]
    foreach ($cellRef in $cellMap.Keys) {
        $ws.Range($cellRef).Value = $cellMap[$cellRef]
    }
}
